$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: update the Lab 4 link text/url
$ws.Range("F9").Value = "[Lab 4 - ER Model](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2025&branch=main&urlpath=tree%2Fdemog180-fa2025%2Flab%2Flab04%2Flab4_clustering_coefficient_er.ipynb)"

# Row 11: remove the "? Problem Set 01" note (moved down to row 12)
$ws.Range("G11").Value = ""

# Row 12: add the "? Problem Set 01" note
$ws.Range("G12").Value = "? Problem Set 01"

# Row 13: topic text change
$ws.Range("D13").Value = "More models: configuration model and stochastic block model"

# Row 14: topic text change
$ws.Range("D14").Value = "Community detection"

# Row 15: topic text change, remove Lab 6 link (moved down to row 17)
$ws.Range("D15").Value = "Midterm review"
$ws.Range("F15").Value = ""

# Row 16: topic text change
$ws.Range("D16").Value = "Midterm"

# Row 17: topic text change, add Lab 6 link
$ws.Range("D17").Value = "Empirical studies of network structure"
$ws.Range("F17").Value = "[Lab 6: Community detection](https://datahub.berkeley.edu/hub/user-redirect/git-pull?repo=https%3A%2F%2Fgithub.com%2Fdfeehan%2Fdemog180-fa2024&branch=main&urlpath=tree%2Fdemog180-fa2024%2Flabs%2Flab06%2Flab06_community_detection.ipynb)"

$wb.Save()
